$wb = $excel.ActiveWorkbook

# --- Competitions sheet: append first six months of 2022 styles -----------
$ws = $wb.Worksheets.Item("Competitions")
$lo = $ws.ListObjects.Item(1)

$newRows = @(
    @(44572, "Open/Club Brews", 2015),
    @(44600, "All Non-American Styles (European, Mexican, etc.)", 2015),
    @(44628, "Lagers and Sour/Funky Beers", 2015),
    @(44663, "American and Hoppy Beers", 2015),
    @(44691, "Open (focus on a ""clone"" beer)", 2015),
    @(44726, "Club Brew/Big Brew Day", 2015)
)

foreach ($entry in $newRows) {
    $row = $lo.ListRows.Add()
    $r = $row.Range
    $r.Cells.Item(1, 1).Value = $entry[0]
    $r.Cells.Item(1, 2).Value = $entry[1]
    $r.Cells.Item(1, 4).Value = $entry[2]
}

# --- Entries sheet: it loses the "active tab" flag -------------------------
$ws2 = $wb.Worksheets.Item("Entries")
[void]$ws2.Range("D126").Select()

# Competitions becomes the active/selected sheet, with D50 selected.
[void]$ws.Activate()
[void]$ws.Range("D50").Select()
